{"js": "// Update the division-problem table: replace each original \"A\u00f7B=\" expression\n// with its new value, per the target diff. Using body.search for each\n// original string (all are unique in the document) and replacing the\n// matched range's text in place keeps formatting (font/size) intact.\nconst replacements = [\n  [\"151\u00f79=\", \"905\u00f79=\"],\n  [\"969\u00f79=\", \"830\u00f79=\"],\n  [\"301\u00f72=\", \"319\u00f74=\"],\n  [\"497\u00f78=\", \"355\u00f75=\"],\n  [\"186\u00f79=\", \"894\u00f76=\"],\n  [\"393\u00f78=\", \"505\u00f75=\"],\n  [\"618\u00f74=\", \"862\u00f75=\"],\n  [\"356\u00f75=\", \"660\u00f73=\"],\n  [\"478\u00f76=\", \"428\u00f79=\"],\n  [\"882\u00f78=\", \"360\u00f73=\"],\n  [\"651\u00f72=\", \"699\u00f78=\"],\n  [\"458\u00f75=\", \"633\u00f72=\"],\n  [\"810\u00f77=\", \"491\u00f78=\"],\n  [\"817\u00f75=\", \"815\u00f75=\"],\n  [\"402\u00f79=\", \"993\u00f77=\"],\n  [\"166\u00f79=\", \"245\u00f72=\"],\n  [\"804\u00f78=\", \"814\u00f79=\"],\n  [\"926\u00f77=\", \"406\u00f75=\"],\n  [\"588\u00f74=\", \"642\u00f76=\"],\n  [\"756\u00f72=\", \"697\u00f74=\"],\n  [\"215\u00f73=\", \"787\u00f75=\"],\n  [\"939\u00f79=\", \"819\u00f74=\"],\n  [\"467\u00f78=\", \"897\u00f73=\"],\n  [\"887\u00f78=\", \"894\u00f73=\"],\n  [\"692\u00f73=\", \"254\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem table: replace each original \"A\u00f7B=\" expression\n# with its new value, per the target diff. Word's Find/Replace (wdReplaceAll)\n# is used for each original string (all are unique in the document), which\n# keeps the surrounding run formatting (font/size) intact.\n$d = $word.ActiveDocument\n\n$map = [ordered]@{\n  \"151\u00f79=\" = \"905\u00f79=\";\n  \"969\u00f79=\" = \"830\u00f79=\";\n  \"301\u00f72=\" = \"319\u00f74=\";\n  \"497\u00f78=\" = \"355\u00f75=\";\n  \"186\u00f79=\" = \"894\u00f76=\";\n  \"393\u00f78=\" = \"505\u00f75=\";\n  \"618\u00f74=\" = \"862\u00f75=\";\n  \"356\u00f75=\" = \"660\u00f73=\";\n  \"478\u00f76=\" = \"428\u00f79=\";\n  \"882\u00f78=\" = \"360\u00f73=\";\n  \"651\u00f72=\" = \"699\u00f78=\";\n  \"458\u00f75=\" = \"633\u00f72=\";\n  \"810\u00f77=\" = \"491\u00f78=\";\n  \"817\u00f75=\" = \"815\u00f75=\";\n  \"402\u00f79=\" = \"993\u00f77=\";\n  \"166\u00f79=\" = \"245\u00f72=\";\n  \"804\u00f78=\" = \"814\u00f79=\";\n  \"926\u00f77=\" = \"406\u00f75=\";\n  \"588\u00f74=\" = \"642\u00f76=\";\n  \"756\u00f72=\" = \"697\u00f74=\";\n  \"215\u00f73=\" = \"787\u00f75=\";\n  \"939\u00f79=\" = \"819\u00f74=\";\n  \"467\u00f78=\" = \"897\u00f73=\";\n  \"887\u00f78=\" = \"894\u00f73=\";\n  \"692\u00f73=\" = \"254\u00f74=\";\n}\n\nforeach ($oldText in $map.Keys) {\n  $newText = $map[$oldText]\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  # wdFindContinue = 1 (Wrap), wdReplaceAll = 2\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
